$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 0) The "_GoBack" bookmark is about to move to the point of the newest edit
#    (the answer paragraph added in step 1 below). Word only ever keeps a
#    single "_GoBack" bookmark, so remove the old one *first*, while it is
#    still unambiguous (there is exactly one "_GoBack" bookmark in the
#    document at this point).
# ---------------------------------------------------------------------------
$oldBkIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "ultimo commit del repositorio remoto") {
        $oldBkIdx = $i + 1
        break
    }
}

if ($oldBkIdx -ne $null) {
    $p = $d.Paragraphs.Item($oldBkIdx)
    $p.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($oldBkIdx)
    $rng = $newPara.Range

    $xmlNoBookmark = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00997D8F" w:rsidRDefault="00997D8F" w:rsidP="00997D8F">
  <w:pPr>
    <w:rPr>
      <w:sz w:val="22"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $rng.InsertXML($xmlNoBookmark)

    $dup = $d.Paragraphs.Item($oldBkIdx + 1)
    $dup.Range.Delete()
}

# ---------------------------------------------------------------------------
# 1) Add the answer paragraph ("Sí, ya que el proyecto posee clases...") plus
#    two blank "Prrafodelista" paragraphs right after the MVC question, i.e.
#    right before the two pre-existing blank paragraphs that precede
#    "Parte 2". This is also where the "_GoBack" bookmark now belongs.
# ---------------------------------------------------------------------------
$mvcQuestionPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Considera usted que el proyecto fue implementado siguiendo los principios del patr") {
        $mvcQuestionPara = $i
        break
    }
}

$targetIndex = $mvcQuestionPara + 1
$anchor = $d.Paragraphs.Item($targetIndex)
$anchor.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($targetIndex)
$rng = $newPara.Range

$xmlAnswer = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:ind w:left="1440"/>
    <w:rPr>
      <w:sz w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>S&#237;</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
    </w:rPr>
    <w:t>, ya que el proyecto posee clases que han sido creadas para una tarea espec&#237;fica, las cuales siguen el paradigma del MVC. El GameController es aquel que realiza los cambios a los modelos (SnakeModel, GoldModel, etc), y posteriormente se actualizan en la interfaz (GameView).</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:ind w:left="1440"/>
    <w:rPr>
      <w:sz w:val="20"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:ind w:left="1440"/>
    <w:rPr>
      <w:sz w:val="20"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xmlAnswer)

# ---------------------------------------------------------------------------
# 2) "Responda a las siguientes interrogantes" now starts a (new) page, so a
#    <w:lastRenderedPageBreak/> marker is recorded in front of its text.
#    Rebuild that paragraph with the marker inserted, then drop the old copy.
# ---------------------------------------------------------------------------
$respondaIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^Responda a las siguientes interrogantes") {
        $respondaIdx = $i
        break
    }
}

$p = $d.Paragraphs.Item($respondaIdx)
$p.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($respondaIdx)
$rng = $newPara.Range

$xmlResponda = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="008E11C2" w:rsidRDefault="008E11C2" w:rsidP="00B61F06">
  <w:pPr>
    <w:rPr>
      <w:sz w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Responda a las siguientes interrogantes</w:t>
  </w:r>
  <w:r w:rsidR="00997D8F">
    <w:rPr>
      <w:sz w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve"> (Solo uno del grupo debe subir este archivo al GitHub con las respuestas grupales)</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
    </w:rPr>
    <w:t>:</w:t>
  </w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xmlResponda)

$dup = $d.Paragraphs.Item($respondaIdx + 1)
$dup.Range.Delete()

# ---------------------------------------------------------------------------
# 3) The old page break now falls elsewhere, so the previously recorded
#    <w:lastRenderedPageBreak/> in front of "El GUIView, el GameController..."
#    is stale and gets removed. Rebuild that paragraph without the marker.
# ---------------------------------------------------------------------------
$guiIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^El GUIView, el GameController") {
        $guiIdx = $i
        break
    }
}

$p = $d.Paragraphs.Item($guiIdx)
$p.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($guiIdx)
$rng = $newPara.Range

$xmlGui = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00A73838" w:rsidRDefault="00A73838" w:rsidP="00A73838">
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:rPr>
      <w:sz w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
    </w:rPr>
    <w:t>El GUIView, el GameController, el GameOverException, y el GameView.</w:t>
  </w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xmlGui)

$dup = $d.Paragraphs.Item($guiIdx + 1)
$dup.Range.Delete()

# ---------------------------------------------------------------------------
# 4) The footer's cached PAGE field result is refreshed from "3" to "1".
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$pageChar = $ftr.Range.Characters.Item(1)
if ($pageChar.Text -eq "3") {
    $pageChar.Text = "1"
}

Write-Output "done"
